$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "PANTA NIMA FREDDY ROLAND JUNIOR",
    "MAZA RIOFRIO CINTHIA NATELAHI",
    "PALACIOS PANTA LUIS MIGUEL",
    "CRISANTO CARMEN ROSITA ABIGAIL",
    "SALAZAR VEGA MARIA FERNANDA",
    "PANTA VARONA CANDY ELIZABETH",
    "VEGA ZAPATA JESUS GABRIEL",
    "ELIAS MACHADO JUANA MARGOT",
    "HIDALGO MOSCOL YESSICA JAZMIN",
    "TALLEDO ELIAS ANDREA ALESSANDRA"
)

$values = @(100, 100, 98, 97, 97, 96, 95, 92, 87, 80)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
